$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter new color-code / legend values in the order that reproduces the
# author's original shared-string insertion sequence.
$ws.Range("D15").Value = "ff0000"
$ws.Range("D14").Value = "00ff44"
$ws.Range("B18").Value = "energy"
$ws.Range("D13").Value = "00aaff"
$ws.Range("G13").Value = "颜色加深"
$ws.Range("H20").Value = "深色"
$ws.Range("E20").Value = "50%+50%"

# Numeric / formatted values
$ws.Range("E13").NumberFormat = "0%"
$ws.Range("E13").Value = 0.5
$ws.Range("F13").Value = 128
$ws.Range("D16").Value = 884400
$ws.Range("D17").Value = 777777

# Duplicate the "颜色加深" label onto row 20 as well
$ws.Range("G20").Value = "颜色加深"

# Weight column (value 1) for rows 20-29
$ws.Range("D20").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("D29").Value = 1

# Update the active selection to reflect the saved cursor position
[void]$ws.Range("E30").Select()
